$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet, defined name, and table (query renamed 39 -> 44) ---
$ws.Name = "query (44)"

$dn = $wb.Names.Item(1)
$dn.Name = "query__44"

$lo = $ws.ListObjects.Item(1)
$lo.Name = "Tabela_query__44"

# --- Fix existing shared string: add one extra trailing space to F154 ---
$ws.Range("F154").Value = "A P Campeao Da Estrada Ltda - 2793748000114, erro na maquina sistema da Quality, erro de vip invalido/não cadastrado.  "

# --- Add new rows 169-187 ---
# Row 169
$ws.Range("A168:T168").Copy($ws.Range("A169:T169"))
$ws.Range("A169").Value = "André"
$ws.Range("B169").ClearContents()
$ws.Range("C169").Value = "Ajuste Lat Long"
$ws.Range("D169").Value = 29017407000103
$ws.Range("E169").Value = "Socape Sociedade Caxiense De Pet Ltda"
$ws.Range("F169").Value = "Levantamento e ajuste de 55 CNPJ's para Lat-Long"
$ws.Range("G169").ClearContents()
$ws.Range("H169").Value = "Torre de Expansão"
$ws.Range("I169").Value = "Teams"
$ws.Range("J169").Value = "Ativo"
$ws.Range("K169").Value = 46007
$ws.Range("L169").Value = 46007
$ws.Range("M169").Value = "NA"
$ws.Range("N169").Value = "Torre"
$ws.Range("O169").ClearContents()
$ws.Range("P169").ClearContents()
$ws.Range("Q169").Value = 0
$ws.Range("R169").ClearContents()
$ws.Range("S169").Value = "Item"
$ws.Range("T169").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"

# Row 170
$ws.Range("A168:T168").Copy($ws.Range("A170:T170"))
$ws.Range("A170").Value = "André"
$ws.Range("B170").Value = "Dúvida"
$ws.Range("C170").Value = "OFS Desconto"
$ws.Range("D170").Value = 11111111111111
$ws.Range("E170").ClearContents()
$ws.Range("F170").Value = "O CN George Otávio, entrou em contato para tirar dúvidas sobre a campanha da Oferta é Seu Desconto. Segundo o consultor, o revendedor reclamou porque entendeu que a campanha era obrigatória, confundindo com a obrigatoriedade da tela do CPF. Expliquei que a obrigatoriedade é apenas para o CPF, e não da campanha."
$ws.Range("G170").Value = "SIM"
$ws.Range("H170").Value = "Torre de Expansão"
$ws.Range("I170").Value = "Teams"
$ws.Range("J170").Value = "Receptivo"
$ws.Range("K170").Value = 46007
$ws.Range("L170").Value = 46007
$ws.Range("M170").Value = "G.N. Rodovia Curitiba"
$ws.Range("N170").Value = "CN"
$ws.Range("O170").ClearContents()
$ws.Range("P170").ClearContents()
$ws.Range("Q170").Value = 0
$ws.Range("R170").ClearContents()
$ws.Range("S170").Value = "Item"
$ws.Range("T170").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(170).RowHeight = 58

# Row 171
$ws.Range("A168:T168").Copy($ws.Range("A171:T171"))
$ws.Range("A171").Value = "André"
$ws.Range("B171").Value = "Solicitação"
$ws.Range("C171").Value = "Instabilidade no APP KMV"
$ws.Range("D171").Value = 78901915001803
$ws.Range("E171").ClearContents()
$ws.Range("F171").Value = "Entrei em contato com 5 postos a fim de solicitar evidências do erro na aba de promoções do Nippo."
$ws.Range("G171").Value = "SIM"
$ws.Range("H171").Value = "Consultor de Negócios"
$ws.Range("I171").Value = "Teams"
$ws.Range("J171").Value = "Ativo"
$ws.Range("K171").Value = 46007
$ws.Range("L171").ClearContents()
$ws.Range("M171").ClearContents()
$ws.Range("N171").Value = "Automação"
$ws.Range("O171").ClearContents()
$ws.Range("P171").ClearContents()
$ws.Range("Q171").Value = 0
$ws.Range("R171").ClearContents()
$ws.Range("S171").Value = "Item"
$ws.Range("T171").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(171).RowHeight = 29

# Row 172
$ws.Range("A168:T168").Copy($ws.Range("A172:T172"))
$ws.Range("A172").Value = "Larissa"
$ws.Range("B172").Value = "Solicitação"
$ws.Range("C172").Value = "Tratativas Financeiras/Repasses"
$ws.Range("D172").Value = 32248017000159
$ws.Range("E172").Value = "Auto Posto Aldeia Da Praia Ltda"
$ws.Range("F172").Value = "Auto Posto Aldeia Da Praia Ltda - 32248017000159 | William pediu apoio com informações sobre Penetração e taxa no GRUPO ROTAS "
$ws.Range("G172").Value = "SIM"
$ws.Range("H172").Value = "Torre de Expansão"
$ws.Range("I172").Value = "Teams"
$ws.Range("J172").Value = "Receptivo"
$ws.Range("K172").Value = 46008
$ws.Range("L172").Value = 46008
$ws.Range("M172").Value = "G.N. Urbano Vitoria"
$ws.Range("N172").Value = "CT"
$ws.Range("O172").ClearContents()
$ws.Range("P172").ClearContents()
$ws.Range("Q172").Value = 0
$ws.Range("R172").ClearContents()
$ws.Range("S172").Value = "Item"
$ws.Range("T172").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(172).RowHeight = 29

# Row 173
$ws.Range("A168:T168").Copy($ws.Range("A173:T173"))
$ws.Range("A173").Value = "Roberta"
$ws.Range("B173").Value = "Reclamação"
$ws.Range("C173").Value = "Erro na impressão do SLIP"
$ws.Range("D173").Value = 1002740000120
$ws.Range("E173").Value = "Auto Posto Leste Ltda"
$ws.Range("F173").Value = "Revendedor relata qie o slip não está registrando a forma de pagaento. Informei que é necessário abrir um chamado com o sistema de automação para confirguar a impressão e assim que estviver com o número do chamado me passar para acompanhar. "
$ws.Range("G173").Value = "SIM"
$ws.Range("H173").Value = "Parceiros Conecta"
$ws.Range("I173").Value = "Whatsapp"
$ws.Range("J173").Value = "Receptivo"
$ws.Range("K173").Value = 46008
$ws.Range("L173").Value = 46008
$ws.Range("M173").Value = "G.N. Rodovia Goiania"
$ws.Range("N173").Value = "Revendedor"
$ws.Range("O173").ClearContents()
$ws.Range("P173").ClearContents()
$ws.Range("Q173").Value = 0
$ws.Range("R173").ClearContents()
$ws.Range("S173").Value = "Item"
$ws.Range("T173").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(173).RowHeight = 43.5

# Row 174
$ws.Range("A168:T168").Copy($ws.Range("A174:T174"))
$ws.Range("A174").Value = "Ricardo"
$ws.Range("B174").Value = "Dúvida"
$ws.Range("C174").Value = "KMV"
$ws.Range("D174").Value = 87091997003136
$ws.Range("E174").Value = "Rede Farroupilha Gestao, Varejo E Inova"
$ws.Range("F174").Value = "O revendedor Cristiano entrou em contato para questionar por que o CNPJ 87.091.997/0031-36 não está refletindo na aba 'Meus Dados' do portal KMV. Solicitei apoio do Thiago Alvarenga para a resolução."
$ws.Range("G174").Value = "NÃO"
$ws.Range("H174").Value = "Coord Engenharia Dig Revendedo"
$ws.Range("I174").Value = "Whatsapp"
$ws.Range("J174").Value = "Receptivo"
$ws.Range("K174").Value = 46008
$ws.Range("L174").ClearContents()
$ws.Range("M174").Value = "G.N. Urbano Porto Alegre"
$ws.Range("N174").Value = "Revendedor"
$ws.Range("O174").ClearContents()
$ws.Range("P174").ClearContents()
$ws.Range("Q174").Value = 0
$ws.Range("R174").ClearContents()
$ws.Range("S174").Value = "Item"
$ws.Range("T174").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(174).RowHeight = 43.5

# Row 175
$ws.Range("A168:T168").Copy($ws.Range("A175:T175"))
$ws.Range("A175").Value = "Roberta"
$ws.Range("B175").Value = "Solicitação"
$ws.Range("C175").Value = "Transação cancelada"
$ws.Range("D175").Value = 37063328000752
$ws.Range("E175").Value = "AutoShoping"
$ws.Range("F175").Value = "Almir solicita confirmação de pagamento p"
$ws.Range("G175").Value = "SIM"
$ws.Range("H175").Value = "Torre de Expansão"
$ws.Range("I175").Value = "Whatsapp"
$ws.Range("J175").Value = "Ativo"
$ws.Range("K175").Value = 46008
$ws.Range("L175").Value = 46008
$ws.Range("M175").Value = "G.N. Urbano Goiania"
$ws.Range("N175").Value = "Revendedor"
$ws.Range("O175").ClearContents()
$ws.Range("P175").ClearContents()
$ws.Range("Q175").Value = 0
$ws.Range("R175").ClearContents()
$ws.Range("S175").Value = "Item"
$ws.Range("T175").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"

# Row 176
$ws.Range("A168:T168").Copy($ws.Range("A176:T176"))
$ws.Range("A176").Value = "Roberta"
$ws.Range("B176").Value = "Solicitação"
$ws.Range("C176").Value = "Transação cancelada"
$ws.Range("D176").Value = 37063328001058
$ws.Range("E176").Value = "AutoShoping"
$ws.Range("F176").Value = "Renato solicita confirmação de transação estornada"
$ws.Range("G176").Value = "SIM"
$ws.Range("H176").Value = "Torre de Expansão"
$ws.Range("I176").Value = "Whatsapp"
$ws.Range("J176").Value = "Receptivo"
$ws.Range("K176").Value = 46008
$ws.Range("L176").Value = 46008
$ws.Range("M176").Value = "G.N. Urbano Goiania"
$ws.Range("N176").Value = "Revendedor"
$ws.Range("O176").ClearContents()
$ws.Range("P176").ClearContents()
$ws.Range("Q176").Value = 0
$ws.Range("R176").ClearContents()
$ws.Range("S176").Value = "Item"
$ws.Range("T176").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"

# Row 177
$ws.Range("A168:T168").Copy($ws.Range("A177:T177"))
$ws.Range("A177").Value = "Thiago"
$ws.Range("B177").Value = "Solicitação"
$ws.Range("C177").Value = "Alteração de Domicílio Bancário"
$ws.Range("D177").Value = 10434857000173
$ws.Range("E177").Value = "Auto Posto W Quatro Ltda"
$ws.Range("F177").Value = "Solciitação de Alteraçao de daados bancários  Rafael Pelegrini de Almeida 321.992.518-96 13/10/1983 15 997271310   Auto Posto GP-7 Ltda 58.968.629/0001-00 Avenida São Paulo, 1988 Bairro: Além Ponte. Sorocaba-SP CEP: 18.013-004 postogp7@gmail.com Banco Itaú 341 Agência 3048 Conta corrente 43326-6     Auto Posto W Quatro Ltda 10.434.857/0001-73 Avenida São Paulo, 2269 Bairro: Além Ponte. Sorocaba-SP CEP: 18.013-004 postowquatro@gmail.com Banco Itaú 341 Agência 3048 Conta corrente 42932-2   Em relação aos telefones dos postos, estamos alterando os números devido ao cancelamento do serviço de linhas “metálicas” na região. Ainda não temos os novos números."
$ws.Range("G177").Value = "SIM"
$ws.Range("H177").Value = "Torre de Expansão"
$ws.Range("I177").Value = "Teams"
$ws.Range("J177").Value = "Receptivo"
$ws.Range("K177").Value = 45979
$ws.Range("L177").Value = 46009
$ws.Range("M177").Value = "G.N. Urbano Campinas"
$ws.Range("N177").Value = "CN"
$ws.Range("O177").ClearContents()
$ws.Range("P177").ClearContents()
$ws.Range("Q177").Value = 0
$ws.Range("R177").ClearContents()
$ws.Range("S177").Value = "Item"
$ws.Range("T177").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(177).RowHeight = 116

# Row 178
$ws.Range("A168:T168").Copy($ws.Range("A178:T178"))
$ws.Range("A178").Value = "Thiago"
$ws.Range("B178").Value = "Solicitação"
$ws.Range("C178").Value = "Treinamento Financeiro"
$ws.Range("D178").Value = 10381595000126
$ws.Range("E178").Value = "Auto Posto Plenitude Ltda"
$ws.Range("F178").Value = "Posto pediu apoio para entender processos do portal e taxas "
$ws.Range("G178").Value = "SIM"
$ws.Range("H178").Value = "Torre de Expansão"
$ws.Range("I178").Value = "Whatsapp"
$ws.Range("J178").Value = "Receptivo"
$ws.Range("K178").Value = 46008
$ws.Range("L178").Value = 46009
$ws.Range("M178").Value = "G.N. Urbano Sp Leste"
$ws.Range("N178").Value = "Revendedor"
$ws.Range("O178").ClearContents()
$ws.Range("P178").ClearContents()
$ws.Range("Q178").Value = 0
$ws.Range("R178").ClearContents()
$ws.Range("S178").Value = "Item"
$ws.Range("T178").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"

# Row 179
$ws.Range("A168:T168").Copy($ws.Range("A179:T179"))
$ws.Range("A179").Value = "Larissa"
$ws.Range("B179").Value = "Reclamação"
$ws.Range("C179").Value = "Transação cancelada"
$ws.Range("D179").Value = 1131154000185
$ws.Range("E179").Value = "P De Abast E Serv V Marques Ltda"
$ws.Range("F179").Value = "P DE ABAST E SERV V MARQUES LTDA - 01131154000185 | Chamado enviado por e-mail pela fabiana  1-48856934200.  ABA-204728690"
$ws.Range("G179").Value = "NÃO"
$ws.Range("H179").Value = "Coord. De Controladoria"
$ws.Range("I179").Value = "E-mail"
$ws.Range("J179").Value = "Receptivo"
$ws.Range("K179").Value = 46008
$ws.Range("L179").ClearContents()
$ws.Range("M179").Value = "G.N. Urbano Rio De Janeiro"
$ws.Range("N179").Value = "Revendedor"
$ws.Range("O179").ClearContents()
$ws.Range("P179").ClearContents()
$ws.Range("Q179").Value = 0
$ws.Range("R179").ClearContents()
$ws.Range("S179").Value = "Item"
$ws.Range("T179").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(179).RowHeight = 29

# Row 180
$ws.Range("A168:T168").Copy($ws.Range("A180:T180"))
$ws.Range("A180").Value = "Thiago"
$ws.Range("B180").Value = "Dúvida"
$ws.Range("C180").Value = "B2C - Cadastro"
$ws.Range("D180").Value = 22988259000125
$ws.Range("E180").Value = "Penta I Pit Stop Auto Posto Ltda"
$ws.Range("F180").Value = "Cliente B2C relatou dificuldades para acessar a conta. Foi orientado a contatar a Central, pois ao consultar o CPF verificamos que já existe cadastro no KMV. Antes disso, sugerimos tentar o procedimento de redefinição de senha."
$ws.Range("G180").Value = "SIM"
$ws.Range("H180").Value = "Central Atendimento B2C"
$ws.Range("I180").Value = "Whatsapp"
$ws.Range("J180").Value = "Receptivo"
$ws.Range("K180").Value = 46009
$ws.Range("L180").Value = 46009
$ws.Range("M180").Value = "G.N. Urbano Sp Abc"
$ws.Range("N180").Value = "Torre"
$ws.Range("O180").ClearContents()
$ws.Range("P180").ClearContents()
$ws.Range("Q180").Value = 0
$ws.Range("R180").ClearContents()
$ws.Range("S180").Value = "Item"
$ws.Range("T180").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(180).RowHeight = 43.5

# Row 181
$ws.Range("A168:T168").Copy($ws.Range("A181:T181"))
$ws.Range("A181").Value = "Larissa"
$ws.Range("B181").Value = "Reclamação"
$ws.Range("C181").Value = "Dúvidas sobre conciliação"
$ws.Range("D181").Value = 16831318000154
$ws.Range("E181").Value = "Machado Amorim P De Gas Com Serv Ltda"
$ws.Range("F181").Value = "Machado Amorim P De Gas Com Serv Ltda - 16831318000154 |  Noraldino (+55 38 9745-3434) Revendedor entrou em contato devido, pois alega que não recebeu repasses do KMV, mas analisando as transações financeiras, foi visto que todos os repasses estão negativos devido as vendas em dinheiro. Solicitando apoio do Gustavo, CN, para visita no posto. "
$ws.Range("G181").Value = "NÃO"
$ws.Range("H181").Value = "Consultor de Negócios"
$ws.Range("I181").Value = "Whatsapp"
$ws.Range("J181").Value = "Receptivo"
$ws.Range("K181").Value = 46007
$ws.Range("L181").Value = 46008
$ws.Range("M181").Value = "G.N. Urbano Uberlandia"
$ws.Range("N181").Value = "Revendedor"
$ws.Range("O181").ClearContents()
$ws.Range("P181").ClearContents()
$ws.Range("Q181").Value = 0
$ws.Range("R181").ClearContents()
$ws.Range("S181").Value = "Item"
$ws.Range("T181").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(181).RowHeight = 58

# Row 182
$ws.Range("A168:T168").Copy($ws.Range("A182:T182"))
$ws.Range("A182").Value = "Larissa"
$ws.Range("B182").Value = "Solicitação"
$ws.Range("C182").Value = "B2C - Duplicidade de trans."
$ws.Range("D182").Value = 25625344000117
$ws.Range("E182").Value = "P Aguia Joia J R Ltda"
$ws.Range("F182").Value = "P Aguia Joia J R Ltda - 25625344000117 | Leandro +55 31 9275-8163, pediu atualizações sobre o chamado 1-48872990560 que o cliente dele abriu dia 16/12"
$ws.Range("G182").Value = "SIM"
$ws.Range("H182").Value = "Central Atendimento B2C"
$ws.Range("I182").Value = "Whatsapp"
$ws.Range("J182").Value = "Receptivo"
$ws.Range("K182").Value = 46008
$ws.Range("L182").Value = 46009
$ws.Range("M182").Value = "G.N. Urbano Belo Horizonte"
$ws.Range("N182").Value = "Revendedor"
$ws.Range("O182").ClearContents()
$ws.Range("P182").ClearContents()
$ws.Range("Q182").Value = 0
$ws.Range("R182").ClearContents()
$ws.Range("S182").Value = "Item"
$ws.Range("T182").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(182).RowHeight = 29

# Row 183
$ws.Range("A168:T168").Copy($ws.Range("A183:T183"))
$ws.Range("A183").Value = "Larissa"
$ws.Range("B183").Value = "Reclamação"
$ws.Range("C183").Value = "B2C - Duplicidade de trans."
$ws.Range("D183").Value = 11398725000104
$ws.Range("E183").Value = "POSTO ABAST GALLENA LAGOA LTDA"
$ws.Range("F183").Value = "POSTO ABAST GALLENA LAGOA LTDA - 11398725000104 |  1-48838811484, cliente reclama ao posto gallena que os seus reembolsos não foram realizados"
$ws.Range("G183").Value = "NÃO"
$ws.Range("H183").Value = "Costumer Hero"
$ws.Range("I183").Value = "Whatsapp"
$ws.Range("J183").Value = "Receptivo"
$ws.Range("K183").Value = 45993
$ws.Range("L183").Value = 46002
$ws.Range("M183").Value = "G.N. Urbano Rio De Janeiro"
$ws.Range("N183").Value = "Revendedor"
$ws.Range("O183").ClearContents()
$ws.Range("P183").ClearContents()
$ws.Range("Q183").Value = 0
$ws.Range("R183").ClearContents()
$ws.Range("S183").Value = "Item"
$ws.Range("T183").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(183).RowHeight = 29

# Row 184
$ws.Range("A168:T168").Copy($ws.Range("A184:T184"))
$ws.Range("A184").Value = "Larissa"
$ws.Range("B184").Value = "Reclamação"
$ws.Range("C184").Value = "B2C - Duplicidade de trans."
$ws.Range("D184").Value = 11398725000104
$ws.Range("E184").Value = "POSTO ABAST GALLENA LAGOA LTDA"
$ws.Range("F184").Value = "POSTO ABAST GALLENA LAGOA LTDA - 11398725000104 |  Cliente JOILSON VIEIRA DE OLIVEIRA segue reclamando sobre o estorno de duas transações que não foram estornadas, enviei o comprovante ao JP"
$ws.Range("G184").Value = "SIM"
$ws.Range("H184").Value = "Costumer Hero"
$ws.Range("I184").Value = "Whatsapp"
$ws.Range("J184").Value = "Receptivo"
$ws.Range("K184").Value = 46009
$ws.Range("L184").Value = 46009
$ws.Range("M184").Value = "G.N. Urbano Rio De Janeiro"
$ws.Range("N184").Value = "Revendedor"
$ws.Range("O184").ClearContents()
$ws.Range("P184").ClearContents()
$ws.Range("Q184").Value = 0
$ws.Range("R184").ClearContents()
$ws.Range("S184").Value = "Item"
$ws.Range("T184").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(184).RowHeight = 43.5

# Row 185
$ws.Range("A168:T168").Copy($ws.Range("A185:T185"))
$ws.Range("A185").Value = "André"
$ws.Range("B185").Value = "Solicitação"
$ws.Range("C185").Value = "Ajuste Lat Long"
$ws.Range("D185").Value = 11111111111111
$ws.Range("E185").ClearContents()
$ws.Range("F185").Value = "Conclusão do ajuste de Lat-Long de 108 CNPJ's; Haron já foi notificado por e-mail. Apontei os 8 CNPJ's que não consegui as informações na notificação."
$ws.Range("G185").Value = "NÃO"
$ws.Range("H185").Value = "Torre de Expansão"
$ws.Range("I185").Value = "Teams"
$ws.Range("J185").Value = "Ativo"
$ws.Range("K185").Value = 46009
$ws.Range("L185").Value = 46009
$ws.Range("M185").ClearContents()
$ws.Range("N185").Value = "Torre"
$ws.Range("O185").ClearContents()
$ws.Range("P185").ClearContents()
$ws.Range("Q185").Value = 0
$ws.Range("R185").ClearContents()
$ws.Range("S185").Value = "Item"
$ws.Range("T185").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(185).RowHeight = 29

# Row 186
$ws.Range("A168:T168").Copy($ws.Range("A186:T186"))
$ws.Range("A186").Value = "André"
$ws.Range("B186").Value = "Reclamação"
$ws.Range("C186").Value = "Pró-frotas"
$ws.Range("D186").Value = 27401974000133
$ws.Range("E186").Value = "Auto Posto Amp Ltda"
$ws.Range("F186").Value = "Reunião com a coordenadora Lucinana do auto posto AMP, Nadine, Lucas Dantas e Larissa. O objetivo foi mapear o fluxo operacional do posto para identificar os erros de preço apontados pelo revendedor. Por fim, por causa de incompatibilidade de agendas, decidimos fazer um grupo no whatsapp a fim de objeter as evidências necessárias. "
$ws.Range("G186").Value = "NÃO"
$ws.Range("H186").Value = "KMV Pista"
$ws.Range("I186").Value = "Teams"
$ws.Range("J186").Value = "Ativo"
$ws.Range("K186").Value = 46009
$ws.Range("L186").ClearContents()
$ws.Range("M186").ClearContents()
$ws.Range("N186").ClearContents()
$ws.Range("O186").ClearContents()
$ws.Range("P186").ClearContents()
$ws.Range("Q186").Value = 0
$ws.Range("R186").ClearContents()
$ws.Range("S186").Value = "Item"
$ws.Range("T186").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"
$ws.Rows.Item(186).RowHeight = 58

# Row 187
$ws.Range("A168:T168").Copy($ws.Range("A187:T187"))
$ws.Range("A187").Value = "André"
$ws.Range("B187").Value = "Solicitação"
$ws.Range("C187").Value = "Ajuste Lat Long"
$ws.Range("D187").Value = 11111111111111
$ws.Range("E187").ClearContents()
$ws.Range("F187").Value = "Entrei em contato com os revendedores a fim de obter a localização dos postos."
$ws.Range("G187").Value = "NÃO"
$ws.Range("H187").Value = "Torre de Expansão"
$ws.Range("I187").Value = "Teams"
$ws.Range("J187").Value = "Ativo"
$ws.Range("K187").Value = 46008
$ws.Range("L187").Value = 46008
$ws.Range("M187").ClearContents()
$ws.Range("N187").Value = "Torre"
$ws.Range("O187").ClearContents()
$ws.Range("P187").ClearContents()
$ws.Range("Q187").Value = 0
$ws.Range("R187").ClearContents()
$ws.Range("S187").Value = "Item"
$ws.Range("T187").Value = "personal/roberta_rocha_ipiranga_ipiranga/Lists/Torre de Expanso  Atendimentos"

# --- Resize table to include new rows, update defined name range, and selection ---
$lo.Resize($ws.Range("A1:T187"))
$dn.RefersTo = "='query (44)'!`$A`$1:`$T`$187"
$ws.Range("A1:T187").Select()
